$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move "SAN DIEGO AREA TOTALS" from B2 to A2, and set B2 to "Totals"
$ws.Range("A2").Value2 = $ws.Range("B2").Value2
$ws.Range("B2").Value2 = "Totals"

# Select column A entirely (mimics clicking the column header)
$ws.Columns("A").Select()

# Column A now holds the same long text as column B, so its best-fit
# width grows to match column B's width.
$ws.Columns("A").ColumnWidth = 21.830729166666668
